{"js": "// Fix a handful of English typos in the \"BACKGROUND\" / \"CONTEXT AND\n// MOTIVATIONS\" sections of the document (per the commit message:\n// \"corrected english background and search\").\nconst body = context.document.body;\n\n// List of [oldText, newText] pairs \u2014 each oldText is a unique\n// substring in the document, so matchCase search + replace is safe.\nconst fixes = [\n  [\"intrdocued\", \"introduced\"],\n  [\"organitzation\", \"organization\"],\n  [\"organizati\u00f3n\", \"organization\"],\n  [\"underestands\", \"understand\"],\n  [\"differnet\", \"different\"],\n  [\"opportunitites\", \"opportunities\"],\n  [\"rutiny\", \"routine\"],\n];\n\nfor (const [oldText, newText] of fixes) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix a handful of English typos in the \"BACKGROUND\" / \"CONTEXT AND\n# MOTIVATIONS\" sections of the document (per the commit message:\n# \"corrected english background and search\").\n\n$d = $word.ActiveDocument\n\n# Each pair is [old typo, corrected text]; every typo is a unique\n# substring in the document, so a simple Find/Replace-all is safe.\n$fixes = @(\n    @(\"intrdocued\", \"introduced\"),\n    @(\"organitzation\", \"organization\"),\n    @(\"organizati\u00f3n\", \"organization\"),\n    @(\"underestands\", \"understand\"),\n    @(\"differnet\", \"different\"),\n    @(\"opportunitites\", \"opportunities\"),\n    @(\"rutiny\", \"routine\")\n)\n\nforeach ($pair in $fixes) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
